$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.028.11'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '2.257.25'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.03'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.94'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.81'
$ws.Range("E10").Value = '  +4.78%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.83'
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").Value = '2.607.13'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.51'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").Value = '2.267.72'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.789'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '41.922.47'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("D20").Value = '0.0₃0901'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.94'
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.40'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.92'
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.50'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.37'
$ws.Range("E28").Value = '  +3.77%  '
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.42'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.79'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.19'
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.15'
$ws.Range("E34").Value = '  +2.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0733'
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.97'
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").Value = '1.958.63'
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.49'
$ws.Range("E45").Value = '  -7.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.90'
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.09'
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '71.46'
$ws.Range("E49").Value = '  -2.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.11'
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E51").Value = '  -2.46%  '
